$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddress, $text) {
    $rng = $ws.Range($rangeAddress)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Enterprises density (per 1000 people): 5.8 -> 5.81, 17.8 -> 17.81
Set-TextValue "C11" "5.81"
Set-TextValue "D11" "17.81"

# Employment (% of total): 14.3 -> 14.31, 35.7 -> 35.69
Set-TextValue "B12" "14.31"
Set-TextValue "C12" "35.69"

# Enterprises (% of total): 63 -> 63.02, 30.5 -> 30.53, 93.5 -> 93.55
Set-TextValue "B14" "63.02"
Set-TextValue "C14" "30.53"
Set-TextValue "D14" "93.55"
